$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Move the "and the following products:" paragraph so that it
#    appears right after the categories ID/Name listing paragraph
#    (and before the Chai / Chef Anton's products ID/Name listing).
# ------------------------------------------------------------------
$moveText = "and the following products:"
$i = 1
$sourcePara = $null
while ($i -le $d.Paragraphs.Count) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq $moveText) {
        $sourcePara = $p
        break
    }
    $i = $i + 1
}

if ($sourcePara -ne $null) {
    $sourcePara.Range.Cut()

    # Re-locate the categories listing paragraph (the one right before
    # the "Chai" products listing) after the cut, and paste right before it.
    $j = 1
    $targetStart = -1
    while ($j -le $d.Paragraphs.Count) {
        $p2 = $d.Paragraphs($j)
        if ($p2.Range.Text.StartsWith("ID: 1, Name: Chai")) {
            $targetStart = $p2.Range.Start
            break
        }
        $j = $j + 1
    }

    if ($targetStart -ge 0) {
        $insertPoint = $d.Range($targetStart, $targetStart)
        $insertPoint.Paste()
    }
}

# ------------------------------------------------------------------
# 2. Colour the "Unknown Method" intro paragraphs red.
# ------------------------------------------------------------------
$d.Content.Find.Execute("You are given a project ConsoleApplication1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$redColor = 255

$k = 1
while ($k -le $d.Paragraphs.Count) {
    $p3 = $d.Paragraphs($k)
    $t = $p3.Range.Text
    if ($t.StartsWith("You are given a project ConsoleApplication1")) {
        $p3.Range.Font.Color = $redColor
    }
    if ($t.StartsWith("You may need to debug and examine the project more closely.")) {
        $p3.Range.Font.Color = $redColor
    }
    $k = $k + 1
}

# ------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark to the edit point inside the phrase
#    "...you need to make the code more readable" (Word drops this
#    bookmark at the location of the most recent edit / typing).
# ------------------------------------------------------------------
$m = 1
$targetPara = $null
while ($m -le $d.Paragraphs.Count) {
    $p4 = $d.Paragraphs($m)
    if ($p4.Range.Text.StartsWith("You are given a project ConsoleApplication1")) {
        $targetPara = $p4
        break
    }
    $m = $m + 1
}

if ($targetPara -ne $null) {
    $searchRange = $targetPara.Range.Duplicate
    $found = $searchRange.Find.Execute("you need to mak", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $splitPos = $searchRange.End
        $bmPoint = $d.Range($splitPos, $splitPos)
        $d.Bookmarks.Add("_GoBack", $bmPoint)
    }
}
